# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the per-locale handback-status sheets to reflect the new
# report generation run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-23 05:15:05"
$zhcn.Range("H2").Value = "2016-03-23 05:15:32"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-23 05:15:11"
$dede.Range("H2").Value = "2016-03-23 05:15:39"
